# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - bound to the (only) Slide Master - "Integral" / "Red Violet"
#   ppt/theme/theme2.xml  - bound only to the Notes Master    - "Office Theme" / "Office"
#
# The target edit swaps the two themes' contents: the Slide Master's theme
# should end up carrying the "Office Theme" palette (the colours that
# currently live in theme2.xml) while the Notes-Master-only theme keeps the
# "Integral" colours. The font scheme and format scheme (fills/lines/effects)
# are identical between the two themes already, so the only substantive,
# automatable change is the 12 theme colours on the reachable (Slide Master)
# theme - which is exactly what PowerPoint's Design > Colors picker would do
# when swapping in the built-in "Office" colour scheme.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
